$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column P (the last existing year column) into the
# new column Q for the header/border row (3) through the last data row (14),
# so the new cells inherit the same number format / font / border as their
# column-P neighbours.
$ws.Range("P3:P14").Copy() | Out-Null
$ws.Range("Q3:Q14").PasteSpecial(-4122) | Out-Null

# Row 3 is just the thin separator row above the year header - no value.
$ws.Range("Q3").Value = ""

# New "2020" year column header.
$ws.Range("Q4").Value = 2020

# New data values for the 2020 column, one per region row.
$ws.Range("Q5").Value = 38.6
$ws.Range("Q6").Value = 42.4
$ws.Range("Q7").Value = 53.2
$ws.Range("Q8").Value = 90.6
$ws.Range("Q9").Value = 52.6
$ws.Range("Q10").Value = 24.5
$ws.Range("Q11").Value = 69.1
$ws.Range("Q12").Value = 32.2
$ws.Range("Q13").Value = 19.1
$ws.Range("Q14").Value = 25.2

# Move the active selection the way the author's session ended up.
$ws.Range("R27").Select() | Out-Null
